$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet (magapoke_2026-02-11)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = 'magapoke_2026-02-18'

# Match page margins used by the other weekly ranking sheets (inches, via points)
$ws.PageSetup.LeftMargin = 0.75 * 72
$ws.PageSetup.RightMargin = 0.75 * 72
$ws.PageSetup.TopMargin = 1 * 72
$ws.PageSetup.BottomMargin = 1 * 72
$ws.PageSetup.HeaderMargin = 0.5 * 72
$ws.PageSetup.FooterMargin = 0.5 * 72

# Populate header + 100 ranking rows (rank, title)
$data = @(
    @('rank', 'title'),
    @('1', 'ブルーロック'),
    @('2', 'WIND BREAKER'),
    @('3', '信じていた仲間達にダンジョン奥地で殺されかけたがギフト『無限ガチャ』でレベル9999の仲間達を手に入れて元パーティーメンバーと世界に復讐＆『ざまぁ！』します！'),
    @('4', '東京卍リベンジャーズ'),
    @('5', 'ベイビーステップ'),
    @('6', 'ギルティサークル'),
    @('7', 'ドラハチ'),
    @('8', 'FAIRY TAIL 100 YEARS QUEST'),
    @('9', 'イレギュラーズ'),
    @('10', '島耕作'),
    @('11', '愛妻の裏アカ'),
    @('12', '十字架のろくにん'),
    @('13', '黄昏町プリズナーズ'),
    @('14', '転生貴族、鑑定スキルで成り上がる～弱小領地を受け継いだので、優秀な人材を増やしていたら、最強領地になってた～'),
    @('15', '君が僕らを悪魔と呼んだ頃'),
    @('16', 'ガチアクタ'),
    @('17', '魔女と傭兵'),
    @('18', '転生したら第七王子だったので、気ままに魔術を極めます'),
    @('19', '魔術ギルド総帥～生まれ変わって今更やり直す2度目の学院生活～'),
    @('20', '【爆アド】生まれた直後から最強悪霊と脳内バトルしてたら魔力量が測定可能域を超えてました〜悪憑の子の謙虚な覇道〜'),
    @('21', 'デッドアカウント'),
    @('22', '異世界ウォーキング'),
    @('23', 'ひゃくえむ。'),
    @('24', 'K-9~警視庁公安部公安第9課異能対策係~'),
    @('25', 'ハードワーカー中田'),
    @('26', '限界集落を脱村した錬金術士、都会で"最強"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～'),
    @('27', '蒼く染めろ'),
    @('28', '追放されなかった男　～二度目の人生は土下座から始まりました～'),
    @('29', '幼馴染とはラブコメにならない'),
    @('30', 'ハンドレッドノート－アグリーダック－'),
    @('31', 'せいぶつ部の田辺くん'),
    @('32', '降り積もれ孤独な死よ'),
    @('33', 'ともだちづくり'),
    @('34', 'ハナバス　苔石花江のバスケ論'),
    @('35', '味方が弱すぎて補助魔法に徹していた宮廷魔法師、追放されて最強を目指す'),
    @('36', '辺境の薬師、都でSランク冒険者となる～英雄村の少年がチート薬で無自覚無双〜'),
    @('37', '屋根の下のアルテミス'),
    @('38', '食糧人類-Starving Anonymous-'),
    @('39', 'お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！'),
    @('40', 'グラぱらっ！'),
    @('41', 'アルキメデスの大戦'),
    @('42', 'アオバノバスケ'),
    @('43', 'ペンの夢に紅をさす'),
    @('44', 'さわらないで小手指くん'),
    @('45', '時々ボソッとロシア語でデレる隣のアーリャさん'),
    @('46', '春くらり'),
    @('47', '普通の本はありません！'),
    @('48', '追放された転生王子、『自動製作《オートクラフト》』スキルで領地を爆速で開拓し最強の村を作ってしまう〜最強クラフトスキルで始める、楽々領地開拓スローライフ〜'),
    @('49', 'Destiny Unchain Online 〜吸血鬼少女となって、やがて『赤の魔王』と呼ばれるようになりました〜'),
    @('50', '恋ニ非ズ'),
    @('51', '異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～'),
    @('52', 'ストーカー行為がバレて人生終了男'),
    @('53', 'ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜'),
    @('54', 'となりの黒川さん'),
    @('55', '南海トラフ巨大地震'),
    @('56', 'なれの果ての僕ら'),
    @('57', 'おやすみ ふみさん'),
    @('58', 'ジュミドロ'),
    @('59', '東京卍リベンジャーズ～場地圭介からの手紙～'),
    @('60', '可愛いだけじゃない式守さん'),
    @('61', '五輪の女神さま 〜なでしこ寮のメダルごはん〜'),
    @('62', '剣帝学院の魔眼賢者'),
    @('63', 'いじめるヤバイ奴'),
    @('64', 'Aランクパーティを離脱した俺は、元教え子たちと迷宮深部を目指す。'),
    @('65', '君が監督！'),
    @('66', '念願の悪役令嬢（ラスボス）の身体を手に入れたぞ！'),
    @('67', '黒猫と魔女の教室'),
    @('68', '白鳥運子は31画'),
    @('69', 'GALAXIAS'),
    @('70', 'おくることば'),
    @('71', '復讐の教科書'),
    @('72', 'シャングリラ・フロンティア～クソゲーハンター、神ゲーに挑まんとす～'),
    @('73', 'MYS'),
    @('74', '不遇職【鑑定士】が実は最強だった～奈落で鍛えた最強の【神眼】で無双する～'),
    @('75', '鳴るさんだぁ'),
    @('76', '冰剣の魔術師が世界を統べる〜世界最強の魔術師である少年は、魔術学院に入学する〜'),
    @('77', '皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～'),
    @('78', '劣等人の魔剣使い　スキルボードを駆使して最強に至る'),
    @('79', 'ヒロインは絶望しました。'),
    @('80', '不遇職『鍛冶師』だけど最強です ～気づけば何でも作れるようになっていた男ののんびりスローライフ～'),
    @('81', 'デスティニーラバーズ'),
    @('82', '日本語が話せないロシア人美少女転入生が頼れるのは、多言語マスターの俺1人'),
    @('83', '阿武ノーマル'),
    @('84', '我間乱 ―修羅―'),
    @('85', 'インフルエンサーにストーカーされています'),
    @('86', '母という呪縛 娘という牢獄'),
    @('87', '捨てられた転生賢者～魔物の森で最強の大魔帝国を作り上げる～'),
    @('88', 'この世界がいずれ滅ぶことを、俺だけが知っている～モンスターが現れた世界で、死に戻りレベルアップ～'),
    @('89', '田んぼで拾った女騎士、田舎で俺の嫁だと思われている'),
    @('90', '魁の花巫女'),
    @('91', '奴隷転生～その奴隷、最強の元王子につき～'),
    @('92', 'ハプスブルク家の華麗なる受難'),
    @('93', '怪病医ラムネ'),
    @('94', 'はじめの一歩'),
    @('95', 'ヤンキーＪＫの異常な愛情'),
    @('96', '金田一少年の事件簿外伝 犯人たちの事件簿'),
    @('97', 'ハンドレッドノート－高校生探偵 天命大地－'),
    @('98', '中華一番！極'),
    @('99', '四刀流の最強配信者～やり込んだVRゲームの設定が現実世界に反映されたので、廃止予定だった戦闘職で無双します～'),
    @('100', 'はっちぽっちぱんち')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

Write-Output "Sheet magapoke_2026-02-18 added with $($data.Count) rows"